$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H40").Value = 1584.4736
$ws_ALC.Range("I40").Value = 1396.5
$ws_ALC.Range("J40").Value = 1793.3334
$ws_ALC.Range("K40").Value = 1396.5
$ws_ALC.Range("L40").Value = 1793.3334
$ws_ALC.Range("M40").Value = -1221.5
$ws_ALC.Range("N40").Value = -2143.3334
$ws_ALC.Range("H96").Value = 539.7857
$ws_ALC.Range("I96").Value = 578.7
$ws_ALC.Range("J96").Value = 442.5
$ws_ALC.Range("K96").Value = 1736.1
$ws_ALC.Range("L96").Value = 1327.5
$ws_ALC.Range("M96").Value = -363.1000000000001
$ws_ALC.Range("N96").Value = -4073.5
$ws_ALC.Range("H132").Value = 10005054
$ws_ALC.Range("I132").Value = 13337712
$ws_ALC.Range("J132").Value = 7080
$ws_ALC.Range("K132").Value = 40013136
$ws_ALC.Range("L132").Value = 21240
$ws_ALC.Range("M132").Value = -40010606
$ws_ALC.Range("N132").Value = -26300
$ws_ALC.Range("H138").Value = 4429.3335
$ws_ALC.Range("I138").Value = 2115.3447
$ws_ALC.Range("J138").Value = 6403.0293
$ws_ALC.Range("K138").Value = 6346.034100000001
$ws_ALC.Range("L138").Value = 19209.0879
$ws_ALC.Range("M138").Value = -1206.034100000001
$ws_ALC.Range("N138").Value = -29489.0879

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H74").Value = 1297
$ws_ARM.Range("I74").Value = 913.41174
$ws_ARM.Range("J74").Value = 2927.25
$ws_ARM.Range("K74").Value = 913.41174
$ws_ARM.Range("L74").Value = 2927.25
$ws_ARM.Range("M74").Value = -39.41174000000001
$ws_ARM.Range("N74").Value = -4675.25
$ws_ARM.Range("H77").Value = 1297
$ws_ARM.Range("I77").Value = 913.41174
$ws_ARM.Range("J77").Value = 2927.25
$ws_ARM.Range("K77").Value = 4567.0587
$ws_ARM.Range("L77").Value = 14636.25
$ws_ARM.Range("M77").Value = -199.0586999999996
$ws_ARM.Range("N77").Value = -23372.25
$ws_ARM.Range("H110").Value = 2498.2
$ws_ARM.Range("I110").Value = 667.7778
$ws_ARM.Range("J110").Value = 5243.8335
$ws_ARM.Range("K110").Value = 667.7778
$ws_ARM.Range("L110").Value = 5243.8335
$ws_ARM.Range("M110").Value = 1377.2222
$ws_ARM.Range("N110").Value = -9333.833500000001
$ws_ARM.Range("H132").Value = 2247.9492
$ws_ARM.Range("I132").Value = 1720.0209
$ws_ARM.Range("J132").Value = 4551.636
$ws_ARM.Range("K132").Value = 5160.0627
$ws_ARM.Range("L132").Value = 13654.908
$ws_ARM.Range("M132").Value = -2630.0627
$ws_ARM.Range("N132").Value = -18714.908

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H94").Value = 897.4
$ws_BSM.Range("I94").Value = 735
$ws_BSM.Range("K94").Value = 735
$ws_BSM.Range("M94").Value = -284
$ws_BSM.Range("H105").Value = 1407.2778
$ws_BSM.Range("I105").Value = 1322.8572
$ws_BSM.Range("J105").Value = 1702.75
$ws_BSM.Range("K105").Value = 1322.8572
$ws_BSM.Range("L105").Value = 1702.75
$ws_BSM.Range("M105").Value = 424.1428000000001
$ws_BSM.Range("N105").Value = -5196.75
$ws_BSM.Range("H134").Value = 2571.3064
$ws_BSM.Range("I134").Value = 2561.1086
$ws_BSM.Range("J134").Value = 2600.625
$ws_BSM.Range("K134").Value = 7683.325800000001
$ws_BSM.Range("L134").Value = 7801.875
$ws_BSM.Range("M134").Value = -5148.325800000001
$ws_BSM.Range("N134").Value = -12871.875

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H22").Value = 318
$ws_CRP.Range("I22").Value = 263.5
$ws_CRP.Range("J22").Value = 499.66666
$ws_CRP.Range("K22").Value = 263.5
$ws_CRP.Range("L22").Value = 499.66666
$ws_CRP.Range("M22").Value = 86.5
$ws_CRP.Range("N22").Value = -1199.66666
$ws_CRP.Range("H132").Value = 2314.3428
$ws_CRP.Range("I132").Value = 2085.9048
$ws_CRP.Range("J132").Value = 2657
$ws_CRP.Range("K132").Value = 6257.714399999999
$ws_CRP.Range("L132").Value = 7971
$ws_CRP.Range("M132").Value = -3727.714399999999
$ws_CRP.Range("N132").Value = -13031
$ws_CRP.Range("H134").Value = 17244124
$ws_CRP.Range("I134").Value = 25002634
$ws_CRP.Range("J134").Value = 2990.3333
$ws_CRP.Range("K134").Value = 75007902
$ws_CRP.Range("L134").Value = 8970.999899999999
$ws_CRP.Range("M134").Value = -75005367
$ws_CRP.Range("N134").Value = -14040.9999

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H87").Value = 13998.5
$ws_CUL.Range("H90").Value = 13998.5
$ws_CUL.Range("H131").Value = 1315.1428
$ws_CUL.Range("I131").Value = 3171.25
$ws_CUL.Range("J131").Value = 1075.6451
$ws_CUL.Range("K131").Value = 9513.75
$ws_CUL.Range("L131").Value = 3226.9353
$ws_CUL.Range("M131").Value = -4473.75
$ws_CUL.Range("N131").Value = -13306.9353

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H132").Value = 3144.6038
$ws_GSM.Range("I132").Value = 2824.5
$ws_GSM.Range("K132").Value = 8473.5
$ws_GSM.Range("M132").Value = -5943.5

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H20").Value = 44114.445
$ws_LTW.Range("I20").Value = 9000
$ws_LTW.Range("J20").Value = 61671.668
$ws_LTW.Range("K20").Value = 9000
$ws_LTW.Range("L20").Value = 61671.668
$ws_LTW.Range("M20").Value = -8774
$ws_LTW.Range("N20").Value = -62123.668
$ws_LTW.Range("H61").Value = 2350.158
$ws_LTW.Range("I61").Value = 549.8570999999999
$ws_LTW.Range("J61").Value = 3400.3333
$ws_LTW.Range("K61").Value = 549.8570999999999
$ws_LTW.Range("L61").Value = 3400.3333
$ws_LTW.Range("M61").Value = -347.8570999999999
$ws_LTW.Range("N61").Value = -3804.3333
$ws_LTW.Range("H113").Value = 2350.158
$ws_LTW.Range("I113").Value = 549.8570999999999
$ws_LTW.Range("J113").Value = 3400.3333
$ws_LTW.Range("K113").Value = 549.8570999999999
$ws_LTW.Range("L113").Value = 3400.3333
$ws_LTW.Range("M113").Value = 1620.1429
$ws_LTW.Range("N113").Value = -7740.3333

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H21").Value = 161256.38
$ws_WVR.Range("J21").Value = 72512.75
$ws_WVR.Range("L21").Value = 72512.75
$ws_WVR.Range("N21").Value = -72982.75
$ws_WVR.Range("H24").Value = 40003.332
$ws_WVR.Range("J24").Value = 40003.332
$ws_WVR.Range("L24").Value = 40003.332
$ws_WVR.Range("N24").Value = -40463.332
$ws_WVR.Range("H28").Value = 25742.25
$ws_WVR.Range("I28").Value = 0
$ws_WVR.Range("J28").Value = 25742.25
$ws_WVR.Range("K28").Value = 0
$ws_WVR.Range("L28").Value = 25742.25
$ws_WVR.Range("M28").ClearContents()
$ws_WVR.Range("N28").Value = -26438.25
$ws_WVR.Range("H30").Value = 27717
$ws_WVR.Range("J30").Value = 29169.834
$ws_WVR.Range("L30").Value = 29169.834
$ws_WVR.Range("N30").Value = -29383.834
$ws_WVR.Range("H35").Value = 161256.38
$ws_WVR.Range("J35").Value = 72512.75
$ws_WVR.Range("L35").Value = 72512.75
$ws_WVR.Range("N35").Value = -73092.75
$ws_WVR.Range("H113").Value = 2916.5
$ws_WVR.Range("I113").Value = 2000
$ws_WVR.Range("J113").Value = 3374.75
$ws_WVR.Range("K113").Value = 6000
$ws_WVR.Range("L113").Value = 10124.25
$ws_WVR.Range("M113").Value = -3830
$ws_WVR.Range("N113").Value = -14464.25
$ws_WVR.Range("H132").Value = 17253.828
$ws_WVR.Range("I132").Value = 2611.44
$ws_WVR.Range("J132").Value = 53859.8
$ws_WVR.Range("K132").Value = 7834.32
$ws_WVR.Range("L132").Value = 161579.4
$ws_WVR.Range("M132").Value = -5304.32
$ws_WVR.Range("N132").Value = -166639.4
$ws_WVR.Range("H136").Value = 3111.4783
$ws_WVR.Range("I136").Value = 3406.0908
$ws_WVR.Range("J136").Value = 2841.4167
$ws_WVR.Range("K136").Value = 10218.2724
$ws_WVR.Range("L136").Value = 8524.250100000001
$ws_WVR.Range("M136").Value = -7668.2724
$ws_WVR.Range("N136").Value = -13624.2501
